# fix(mhs): memperbaiki import mhs ditambah gel dan tahun
#
# Adds "GEL" and "TAHUN" columns to the mahasiswa import template (Sheet1),
# clears the stale/blank placeholder formatting left over in rows 2-4 for
# columns A, B, C, F (column G keeps its formatting), and tidies up the
# leftover per-column styling on the two unused helper sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# New header columns: GEL (gelombang) and TAHUN (tahun), right after PRODI.
$ws1.Range("J1").Value = "GEL"
$ws1.Range("K1").Value = "TAHUN"

# Drop the leftover formatted-but-empty placeholder cells in rows 2-4 that
# no longer belong in the cleaned-up template (column G keeps its style).
$ws1.Range("A2:C4").Clear()
$ws1.Range("F2:F4").Clear()

# Match the author's on-screen selection after the edit.
$ws1.Range("J16").Select()

# Sheet2/Sheet3 are empty helper sheets; drop their now-unused per-column
# style overrides so only the plain width customization remains.
$ws2.Cells.ClearFormats()
$ws3.Cells.ClearFormats()
